# mau_bangchamcong_chitiet_chunhat_chot_kx.xlsx
# "fix tải các bảng công CHủ Nhật"
#
# The sheet's header row (row 3) listed several columns that are no longer
# produced by the export (Cấp bậc, Số phút ca, Phút nghỉ phép, Phút tăng ca
# 100%/150%, Phút tăng ca đêm, Phút nghỉ không lương, Phân loại, HC Category)
# and is missing a "Phút tăng ca 200%" column. This script removes the
# obsolete columns and inserts the new one in their place, matching the
# template used for other "chốt" timesheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the obsolete columns, right-to-left so earlier column letters stay
# valid while we work.
#   G  -> Cấp bậc
#   J  -> Số phút ca
#   N  -> Phút nghỉ phép
#   O  -> Phút tăng ca 100%
#   P  -> Phút tăng ca 150%
#   Q  -> Phút tăng ca đêm
#   R  -> Phút nghỉ không lương
#   U  -> Phân loại
#   V  -> HC Category
$colsToDelete = @("V", "U", "R", "Q", "P", "O", "N", "J", "G")
foreach ($col in $colsToDelete) {
    $ws.Range($col + "1").EntireColumn.Delete()
}

# After the deletions the header row reads:
#   ... F=Phòng ban, G=Ngày, H=Ca, I=Giờ vào, J=Giờ ra, K=Phút hành chính,
#   L=Phút nghỉ khác, M=Loại nghỉ khác
# Insert a new column before the old L so we can add "Phút tăng ca 200%"
# right after "Phút hành chính".
$ws.Range("L1").EntireColumn.Insert()
$ws.Range("L3").Value = "Phút tăng ca 200%"
# Give the new column the same display width (~19 characters, like the
# "Phút tăng ca 100%" column it replaces) used elsewhere in this template.
$ws.Range("L1").EntireColumn.ColumnWidth = 18.16666666667

# Mirror the interactive selection state left behind in the saved file
# (a full-column selection on column I).
$ws.Range("I3").EntireColumn.Select()
